# Reporte semanal semana 2 ciclo 1
# Fill in the "Total de Horas Trabajadas" value for the Quality/Process
# Manager row on the second weekly-log sheet ("logt-2"), which was left
# blank in the previous report.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("logt-2")
$ws.Activate() | Out-Null

$ws.Range("C4").Value = 6

# Leave the selection where the user would land after typing the value
# and pressing Enter.
$ws.Range("C5").Select() | Out-Null
